$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Amira Sobhy, Dr. Veronia Rafat, Administrator, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel'
$ws.Range("G3").Value = 'Dr. Majorelle Magdy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Administrator, Dr. Asmaa Reda'
$ws.Range("G4").Value = 'Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel, Dr. Asmaa Reda'
$ws.Range("G5").Value = 'Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Amira Sobhy'
$ws.Range("G6").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Manar Montaser'
$ws.Range("G7").Value = 'Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Kerelos Zareef'
$ws.Range("G8").Value = 'Dr. Nada Mohammad, Dr. Abeer Ragab'
$ws.Range("G11").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G12").Value = 'Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Eman M. Abo-Sakaya'
$ws.Range("G13").Value = 'Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa'
$ws.Range("G15").Value = 'Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat'
$ws.Range("G17").Value = 'Dr. Esraa Samy, Dr. Mohammad Safwat'
$ws.Range("G19").Value = 'Dr. Rania Ahmad Youssef, Dr. Mariam Toma Gerges'
$ws.Range("G20").Value = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Range("G27").Value = 'Dr. Hana Amr, Dr. Nourham Mostafa'
$ws.Range("G28").Value = 'Dr. Maryam Ashraf, Dr. Aya Emad'
$ws.Range("G30").Value = 'Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Shorok Mohammad'
